# Auto-generated Excel COM-interop script
# Applies market-data value updates to the Coeurl_Profits workbook (per sheet: ALC, ARM, BSM, CRP, CUL, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 362.22223
$ws.Range("I55").Value = 376.66666
$ws.Range("J55").Value = 333.33334
$ws.Range("K55").Value = 376.66666
$ws.Range("L55").Value = 333.33334
$ws.Range("M55").Value = -162.66666
$ws.Range("N55").Value = -761.33334

$ws.Range("H80").Value = 1143.7273
$ws.Range("I80").Value = 893.5
$ws.Range("K80").Value = 2680.5
$ws.Range("M80").Value = -1682.5

$ws.Range("H83").Value = 1143.7273
$ws.Range("I83").Value = 893.5
$ws.Range("K83").Value = 8041.5
$ws.Range("M83").Value = -3049.5

$ws.Range("H94").Value = 39303.918
$ws.Range("I94").Value = 39303.918
$ws.Range("K94").Value = 39303.918
$ws.Range("M94").Value = -38852.918

$ws.Range("H131").Value = 9684.637000000001
$ws.Range("J131").Value = 80002.5
$ws.Range("L131").Value = 240007.5
$ws.Range("N131").Value = -250087.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33336598
$ws.Range("I2").Value = 40003252
$ws.Range("J2").Value = 3318
$ws.Range("K2").Value = 40003252
$ws.Range("L2").Value = 3318
$ws.Range("M2").Value = -40003139
$ws.Range("N2").Value = -3544

$ws.Range("H24").Value = 19785
$ws.Range("J24").Value = 19785
$ws.Range("L24").Value = 19785
$ws.Range("N24").Value = -20533

$ws.Range("H32").Value = 4567.522
$ws.Range("I32").Value = 4193.4185
$ws.Range("K32").Value = 4193.4185
$ws.Range("M32").Value = -3906.4185

$ws.Range("H63").Value = 4744.143
$ws.Range("I63").Value = 4744.143
$ws.Range("K63").Value = 4744.143
$ws.Range("M63").Value = -4058.143

$ws.Range("H66").Value = 4744.143
$ws.Range("I66").Value = 4744.143
$ws.Range("K66").Value = 23720.715
$ws.Range("M66").Value = -20288.715

$ws.Range("H74").Value = 7387.871
$ws.Range("I74").Value = 1597.7307
$ws.Range("K74").Value = 1597.7307
$ws.Range("M74").Value = -723.7307000000001

$ws.Range("H77").Value = 7387.871
$ws.Range("I77").Value = 1597.7307
$ws.Range("K77").Value = 7988.6535
$ws.Range("M77").Value = -3620.6535

$ws.Range("H100").Value = 19785
$ws.Range("J100").Value = 19785
$ws.Range("L100").Value = 19785
$ws.Range("N100").Value = -21949

$ws.Range("H102").Value = 2531.5186
$ws.Range("I102").Value = 1841.4348
$ws.Range("J102").Value = 6499.5
$ws.Range("K102").Value = 1841.4348
$ws.Range("L102").Value = 6499.5
$ws.Range("M102").Value = -219.4348
$ws.Range("N102").Value = -9743.5

$ws.Range("H116").Value = 33336598
$ws.Range("I116").Value = 40003252
$ws.Range("J116").Value = 3318
$ws.Range("K116").Value = 40003252
$ws.Range("L116").Value = 3318
$ws.Range("M116").Value = -40000958
$ws.Range("N116").Value = -7906

$ws.Range("H122").Value = 2803.6667
$ws.Range("I122").Value = 2513.2778
$ws.Range("K122").Value = 7539.8334
$ws.Range("M122").Value = -5089.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33336598
$ws.Range("I3").Value = 40003252
$ws.Range("J3").Value = 3318
$ws.Range("K3").Value = 40003252
$ws.Range("L3").Value = 3318
$ws.Range("M3").Value = -40003138
$ws.Range("N3").Value = -3546

$ws.Range("H43").Value = 500684
$ws.Range("J43").Value = 500684
$ws.Range("L43").Value = 500684
$ws.Range("N43").Value = -501046

$ws.Range("H86").Value = 1897.5
$ws.Range("I86").Value = 1897.5
$ws.Range("K86").Value = 1897.5
$ws.Range("M86").Value = -774.5

$ws.Range("H89").Value = 1897.5
$ws.Range("I89").Value = 1897.5
$ws.Range("K89").Value = 9487.5
$ws.Range("M89").Value = -3871.5

$ws.Range("H99").Value = 4138.55
$ws.Range("I99").Value = 1303
$ws.Range("K99").Value = 1303
$ws.Range("M99").Value = 195

$ws.Range("H107").Value = 2629.3076
$ws.Range("I107").Value = 2459.2
$ws.Range("J107").Value = 2735.625
$ws.Range("K107").Value = 2459.2
$ws.Range("L107").Value = 2735.625
$ws.Range("M107").Value = -539.1999999999998
$ws.Range("N107").Value = -6575.625

$ws.Range("H134").Value = 2699
$ws.Range("I134").Value = 2176.7036
$ws.Range("K134").Value = 6530.110799999999
$ws.Range("M134").Value = -3995.110799999999

$ws.Range("H140").Value = 97990
$ws.Range("J140").Value = 97990
$ws.Range("L140").Value = 97990
$ws.Range("N140").Value = -108350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2832.3142
$ws.Range("I31").Value = 2066.1304
$ws.Range("K31").Value = 2066.1304
$ws.Range("M31").Value = -1771.1304

$ws.Range("H34").Value = 2832.3142
$ws.Range("I34").Value = 2066.1304
$ws.Range("K34").Value = 2066.1304
$ws.Range("M34").Value = -1864.1304

$ws.Range("H62").Value = 3509.8
$ws.Range("I62").Value = 3633
$ws.Range("K62").Value = 3633
$ws.Range("M62").Value = -3009

$ws.Range("H65").Value = 3509.8
$ws.Range("I65").Value = 3633
$ws.Range("K65").Value = 18165
$ws.Range("M65").Value = -15045

$ws.Range("H99").Value = 6127.5
$ws.Range("I99").Value = 2804
$ws.Range("K99").Value = 2804
$ws.Range("M99").Value = -1306

$ws.Range("H122").Value = 927.8461
$ws.Range("I122").Value = 996.13635
$ws.Range("J122").Value = 552.25
$ws.Range("K122").Value = 2988.40905
$ws.Range("L122").Value = 1656.75
$ws.Range("M122").Value = -538.4090500000002
$ws.Range("N122").Value = -6556.75

$ws.Range("H126").Value = 6127.5
$ws.Range("I126").Value = 2804
$ws.Range("K126").Value = 8412
$ws.Range("M126").Value = -5942

$ws.Range("H134").Value = 10097.777
$ws.Range("I134").Value = 4467.282
$ws.Range("K134").Value = 13401.846
$ws.Range("M134").Value = -10866.846

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3968917.5
$ws.Range("I34").Value = 94.5
$ws.Range("J34").Value = 4630388
$ws.Range("K34").Value = 283.5
$ws.Range("L34").Value = 13891164
$ws.Range("M34").Value = -199.5
$ws.Range("N34").Value = -13891332

$ws.Range("H39").Value = 1699.375
$ws.Range("J39").Value = 1699.375
$ws.Range("L39").Value = 5098.125
$ws.Range("N39").Value = -5686.125

$ws.Range("H44").Value = 762.3
$ws.Range("I44").Value = 310.6
$ws.Range("K44").Value = 931.8000000000001
$ws.Range("M44").Value = -533.8000000000001

$ws.Range("H55").Value = 1123.5714
$ws.Range("J55").Value = 2129.4285
$ws.Range("L55").Value = 6388.2855
$ws.Range("N55").Value = -6742.2855

$ws.Range("H70").Value = 2750
$ws.Range("I70").Value = 2750
$ws.Range("K70").Value = 8250
$ws.Range("M70").Value = -7935

$ws.Range("H73").Value = 2750
$ws.Range("I73").Value = 2750
$ws.Range("K73").Value = 8250
$ws.Range("M73").Value = -7158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5472.4546
$ws.Range("I40").Value = 4434.7646
$ws.Range("K40").Value = 4434.7646
$ws.Range("M40").Value = -4298.7646

$ws.Range("H68").Value = 3000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 3000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H93").Value = 2032.1111
$ws.Range("I93").Value = 2215
$ws.Range("J93").Value = 1392
$ws.Range("K93").Value = 2215
$ws.Range("L93").Value = 1392
$ws.Range("M93").Value = -967
$ws.Range("N93").Value = -3888

$ws.Range("H122").Value = 5466.5835
$ws.Range("I122").Value = 4641.294
$ws.Range("J122").Value = 7470.857
$ws.Range("K122").Value = 13923.882
$ws.Range("L122").Value = 22412.571
$ws.Range("M122").Value = -11473.882
$ws.Range("N122").Value = -27312.571

$ws.Range("H133").Value = 59879.8
$ws.Range("J133").Value = 59879.8
$ws.Range("L133").Value = 59879.8
$ws.Range("N133").Value = -64939.8

$ws.Range("H136").Value = 8125.375
$ws.Range("I136").Value = 7667.1665
$ws.Range("K136").Value = 23001.4995
$ws.Range("M136").Value = -20451.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 694.5
$ws.Range("I107").Value = 765.1111
$ws.Range("K107").Value = 2295.3333
$ws.Range("M107").Value = -375.3332999999998

$ws.Range("H125").Value = 250023740
$ws.Range("J125").Value = 250023740
$ws.Range("L125").Value = 250023740
$ws.Range("N125").Value = -250033580
